$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new race-result row (row 14)
$ws.Range("A13").Copy($ws.Range("A14")) | Out-Null
$ws.Range("A14").Value = "2/16/2026"
$ws.Range("B14").Value = "UAE Tour"
$ws.Range("C14").Value = "Stage 1"
$ws.Range("D14").Value = "Isaac Del Toro"
$ws.Range("E14").Value = "Cees Bol"
$ws.Range("F14").Value = "Antonio Tiberi"
$ws.Range("G14").Value = "Matteo Malucelli"
$ws.Range("H14").Value = "Alessandro Romele"
$ws.Range("I14").Value = "Ilan Van Wilder"
$ws.Range("J14").Value = "Riley Pickrell"
$ws.Range("K14").Value = "Lennert Van Eetvelt"
$ws.Range("L14").Value = "Ethan Vernon"
$ws.Range("M14").Value = "Tobias Halland Johannessen"

$ws.Range("D20").Select() | Out-Null
